$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("F2").Value = 4.4
    $ws.Range("G2").Value = 5.7
    $ws.Range("H2").Value = 1.78
    $ws.Range("I2").Value = 2.02
    $ws.Range("J2").Value = 3.5
    $ws.Range("L2").Value = 1.32
    $ws.Range("P2").Value = 2.08
    $ws.Range("V2").Value = 1.99
    $ws.Range("W2").Value = 1.24
    $ws.Range("F3").Value = 3
    $ws.Range("I3").Value = 3.25
    $ws.Range("J3").Value = 2.74
    $ws.Range("K3").Value = 2.92
    $ws.Range("L3").Value = 1.61
    $ws.Range("T3").Value = 2.34
    $ws.Range("U3").Value = 1.66
    $ws.Range("Y3").Value = 7.8
    $ws.Range("AD3").Value = 1000
    $ws.Range("AO3").Value = 100
    $ws.Range("F4").Value = 1.71
    $ws.Range("H4").Value = 4.8
    $ws.Range("I4").Value = 5.8
    $ws.Range("J4").Value = 3.7
    $ws.Range("L4").Value = 1.33
    $ws.Range("N4").Value = 3.8
    $ws.Range("P4").Value = 1.98
    $ws.Range("Q4").Value = 1.77
    $ws.Range("R4").Value = 1.38
    $ws.Range("S4").Value = 3.1
    $ws.Range("T4").Value = 1.78
    $ws.Range("U4").Value = 2.02
    $ws.Range("V4").Value = 1.21
    $ws.Range("W4").Value = 2.2
    $ws.Range("Y4").Value = 19.5
    $ws.Range("Z4").Value = 1000
    $ws.Range("AC4").Value = 9.4
    $ws.Range("AD4").Value = 22
    $ws.Range("AH4").Value = 21
    $ws.Range("AI4").Value = 1000
    $ws.Range("AN4").Value = 11.5
    $ws.Range("G5").Value = 1.77
    $ws.Range("K5").Value = 4.2
    $ws.Range("L5").Value = 1.37
    $ws.Range("Q5").Value = 2.06
    $ws.Range("W5").Value = 2.28
    $ws.Range("I6").Value = 17
    $ws.Range("L6").Value = 1.23
    $ws.Range("S6").Value = 2.06
    $ws.Range("V6").Value = 1.06
    $ws.Range("X6").Value = 1000
    $ws.Range("Y6").Value = 1000
    $ws.Range("AG6").Value = 13.5
    $ws.Range("AH6").Value = 1000
    $ws.Range("AL6").Value = 1000
    $ws.Range("AN6").Value = 4.1
    $ws.Range("J7").Value = 8.4
    $ws.Range("AC7").Value = 24
    $ws.Range("AO7").Value = 3.65
    $ws.Range("G8").Value = 3
    $ws.Range("K8").Value = 2.96
    $ws.Range("L8").Value = 1.69
    $ws.Range("M8").Value = 1.17
    $ws.Range("W8").Value = 1.5
    $ws.Range("Z8").Value = 25
    $ws.Range("AK8").Value = 55
    $ws.Range("I9").Value = 7.4
    $ws.Range("S9").Value = 4.4
    $ws.Range("AH9").Value = 1000
    $ws.Range("AJ9").Value = 24
    $ws.Range("AN9").Value = 21
    $ws.Range("G10").Value = 1.79
    $ws.Range("I10").Value = 6
    $ws.Range("K10").Value = 4.1
    $ws.Range("L10").Value = 1.41
    $ws.Range("N10").Value = 3.65
    $ws.Range("O10").Value = 1.32
    $ws.Range("P10").Value = 1.91
    $ws.Range("R10").Value = 1.34
    $ws.Range("U10").Value = 1.96
    $ws.Range("V10").Value = 1.2
    $ws.Range("W10").Value = 2.26
    $ws.Range("Z10").Value = 50
    $ws.Range("AA10").Value = 160
    $ws.Range("AF10").Value = 11
    $ws.Range("AH10").Value = 970
    $ws.Range("AI10").Value = 85
    $ws.Range("AM10").Value = 140
    $ws.Range("AN10").Value = 13.5
    $ws.Range("F11").Value = 2.5
    $ws.Range("G11").Value = 2.84
    $ws.Range("I11").Value = 3.7
    $ws.Range("M11").Value = 1.13
    $ws.Range("S11").Value = 5.5
    $ws.Range("X11").Value = 9.800000000000001
    $ws.Range("AE11").Value = 60
    $ws.Range("AL11").Value = 75
    $ws.Range("F12").Value = 2.22
    $ws.Range("G12").Value = 2.3
    $ws.Range("K12").Value = 3.5
    $ws.Range("P12").Value = 1.83
    $ws.Range("Q12").Value = 2.12
    $ws.Range("U12").Value = 2.08
    $ws.Range("W12").Value = 1.76
    $ws.Range("AH12").Value = 23
    $ws.Range("AN12").Value = 26
    $ws.Range("AO12").Value = 1000
    $ws.Range("G13").Value = 2.6
    $ws.Range("J13").Value = 2.72
    $ws.Range("L13").Value = 1.5
    $ws.Range("V13").Value = 1.28
    $ws.Range("W13").Value = 1.63
    $ws.Range("N14").Value = 3.6
    $ws.Range("O14").Value = 1.35
    $ws.Range("P14").Value = 1.89
    $ws.Range("Q14").Value = 2.02
    $ws.Range("S14").Value = 3.7
    $ws.Range("T14").Value = 2.02
    $ws.Range("U14").Value = 1.9
    $ws.Range("W14").Value = 2.42
    $ws.Range("AC14").Value = 8.800000000000001
    $ws.Range("AM14").Value = 150
    $ws.Range("G15").Value = 1.51
    $ws.Range("I15").Value = 10
    $ws.Range("K15").Value = 4.7
    $ws.Range("N15").Value = 3.35
    $ws.Range("O15").Value = 1.37
    $ws.Range("W15").Value = 2.96
    $ws.Range("AI15").Value = 200
    $ws.Range("AO15").Value = 420
    $ws.Range("K16").Value = 4.4
    $ws.Range("L16").Value = 1.47
    $ws.Range("N16").Value = 2.96
    $ws.Range("O16").Value = 1.45
    $ws.Range("P16").Value = 1.69
    $ws.Range("Q16").Value = 2.32
    $ws.Range("R16").Value = 1.23
    $ws.Range("S16").Value = 4.7
    $ws.Range("T16").Value = 2.4
    $ws.Range("U16").Value = 1.58
    $ws.Range("X16").Value = 12.5
    $ws.Range("Y16").Value = 27
    $ws.Range("AA16").Value = 590
    $ws.Range("AB16").Value = 6.2
    $ws.Range("AC16").Value = 12
    $ws.Range("AD16").Value = 48
    $ws.Range("AE16").Value = 290
    $ws.Range("AF16").Value = 7.6
    $ws.Range("AH16").Value = 970
    $ws.Range("AI16").Value = 270
    $ws.Range("AL16").Value = 65
    $ws.Range("AM16").Value = 380
